$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update relay settings per October 2016 demo changes
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 15
$ws.Range("G2").Value = 0.8

$ws.Range("D8").Value = 9000
$ws.Range("D9").Value = 10000

# Update the saved cursor/selection position
$ws.Range("H2").Select()
